$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Ui test cases"
$ws1.Range("D2").Value = "ADV_UI_1.0"
$ws1.Range("G2").Value = "Advertisment Ui test suite"
$ws1.Range("G6").Value = 44795

$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Api test cases"

$ws2.Range("D2").Value = "ADV_API_1.0"
$ws2.Range("G2").Value = "Advertisment Api test Suite"

$ws2.Range("C13").Value = "User is able to add advertisement with POST request"
$ws2.Range("C14").Value = "User is able to get list of all advertisement with GET request"
$ws2.Range("C15").Value = "User is able to edit advertisement with PATCH request"
$ws2.Range("C16").Value = "User is able to edit advertisement with PUT request"
$ws2.Range("C17").Value = "User is able to get info about specific advertisement with GET request"
$ws2.Range("C18").Value = "Error is shown in response if mandatory fields are missing from request body for POST request"
$ws2.Range("C19").Value = "Error code `"404`" is shown if incorrect request url is send for any method`n"
$ws2.Range("C20").Value = "Error is not shown if user sends a put request without whole body.`nMandatory parameter like name and price"
$ws2.Range("C21").Value = "Error is not shown if user sends a put request without whole body.`nMandatory parameter like name and price"
$ws2.Range("C22").Value = "Error is shown if user passes _id in body for put request"

Write-Output "done"
$ws2.Range("B23:G28").ClearContents()
Write-Output "cleared"
$ws2.Range("C23:D28").UnMerge()
$ws2.Range("E23:G28").UnMerge()
Write-Output "unmerged"
